$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 10.77061114262485
$ws.Range("C2").Value = 6.491500113558428
$ws.Range("D2").Value = 14.70524273902173
$ws.Range("E2").Value = 15.88420249830786
$ws.Range("G2").Value = 3.686658768494158
$ws.Range("I2").Value = 25.59891293014886
$ws.Range("J2").Value = 9.270518408496395
$ws.Range("K2").Value = 11.39595195844102
$ws.Range("M2").Value = 17.39578954905783
$ws.Range("N2").Value = 20.85964984281771
$ws.Range("O2").Value = 28.74842398291112
# Row 3
$ws.Range("B3").Value = 10.55808891684096
$ws.Range("C3").Value = 6.362660701826678
$ws.Range("D3").Value = 14.67935236086104
$ws.Range("E3").Value = 15.88285012322617
$ws.Range("G3").Value = 3.688631850424942
$ws.Range("I3").Value = 25.67605024841424
$ws.Range("J3").Value = 9.287631143554622
$ws.Range("K3").Value = 11.26544493208846
$ws.Range("M3").Value = 17.34965333848378
$ws.Range("N3").Value = 20.92121885688772
$ws.Range("O3").Value = 28.81534326509517
# Row 4
$ws.Range("B4").Value = 10.42760511742554
$ws.Range("C4").Value = 6.283582718600288
$ws.Range("D4").Value = 14.66637632591465
$ws.Range("E4").Value = 15.88488431456742
$ws.Range("G4").Value = 3.68990823961268
$ws.Range("I4").Value = 25.72760975483686
$ws.Range("J4").Value = 9.299024161221709
$ws.Range("K4").Value = 11.18669920640082
$ws.Range("M4").Value = 17.32425271090301
$ws.Range("N4").Value = 20.96079056794346
$ws.Range("O4").Value = 28.86165273010802
# Row 5
$ws.Range("B5").Value = 10.37450951454562
$ws.Range("C5").Value = 6.251411100493143
$ws.Range("D5").Value = 14.66182711587859
$ws.Range("E5").Value = 15.88643484845703
$ws.Range("G5").Value = 3.690444751818001
$ws.Range("I5").Value = 25.74967531040014
$ws.Range("J5").Value = 9.303889968323602
$ws.Range("K5").Value = 11.15499334916076
$ws.Range("M5").Value = 17.31464543964606
$ws.Range("N5").Value = 20.97736224682296
$ws.Range("O5").Value = 28.88183538166786
# Row 6
$ws.Range("B6").Value = 10.36569997135456
$ws.Range("C6").Value = 6.24607360096266
$ws.Range("D6").Value = 14.6611164450693
$ws.Range("E6").Value = 15.88673592612539
$ws.Range("G6").Value = 3.69053482967709
$ws.Range("I6").Value = 25.75340295582404
$ws.Range("J6").Value = 9.304711413998721
$ws.Range("K6").Value = 11.14975280102374
$ws.Range("M6").Value = 17.31309528666443
$ws.Range("N6").Value = 20.98014092967426
$ws.Range("O6").Value = 28.88526583645242
# Row 7
$ws.Range("B7").Value = 10.42688863563172
$ws.Range("C7").Value = 6.283148563999238
$ws.Range("D7").Value = 14.66631197790363
$ws.Range("E7").Value = 15.88490230250134
$ws.Range("G7").Value = 3.689915408840617
$ws.Range("I7").Value = 25.72790306940205
$ws.Range("J7").Value = 9.299088879591871
$ws.Range("K7").Value = 11.18627001078926
$ws.Range("M7").Value = 17.32412012312328
$ws.Range("N7").Value = 20.96101225233042
$ws.Range("O7").Value = 28.86191961381764
# Row 8
$ws.Range("B8").Value = 10.69737749977981
$ws.Range("C8").Value = 6.447097247976798
$ws.Range("D8").Value = 14.69571221224784
$ws.Range("E8").Value = 15.8831425134159
$ws.Range("G8").Value = 3.687325646272862
$ws.Range("I8").Value = 25.6246383997355
$ws.Range("J8").Value = 9.27623524890018
$ws.Range("K8").Value = 11.35068670382924
$ws.Range("M8").Value = 17.37927900945259
$ws.Range("N8").Value = 20.88051264975856
$ws.Range("O8").Value = 28.77041324201818
# Row 9
$ws.Range("B9").Value = 11.22459534427381
$ws.Range("C9").Value = 6.766879830721533
$ws.Range("D9").Value = 14.77633271334669
$ws.Range("E9").Value = 15.90235302431451
$ws.Range("G9").Value = 3.682759809550784
$ws.Range("I9").Value = 25.4554680369247
$ws.Range("J9").Value = 9.238432327675218
$ws.Range("K9").Value = 11.682536965946
$ws.Range("M9").Value = 17.51030930889708
$ws.Range("N9").Value = 20.73662242164059
$ws.Range("O9").Value = 28.63246274100099
# Row 10
$ws.Range("B10").Value = 11.60579532029983
$ws.Range("C10").Value = 6.998253648597574
$ws.Range("D10").Value = 14.8492513654799
$ws.Range("E10").Value = 15.93017404639568
$ws.Range("G10").Value = 3.679714554305359
$ws.Range("I10").Value = 25.35153995807725
$ws.Range("J10").Value = 9.214913628249882
$ws.Range("K10").Value = 11.92997593181438
$ws.Range("M10").Value = 17.6199922373697
$ws.Range("N10").Value = 20.63933715189208
$ws.Range("O10").Value = 28.55649372242025
# Row 11
$ws.Range("B11").Value = 11.77705087358459
$ws.Range("C11").Value = 7.102237929015372
$ws.Range("D11").Value = 14.88531704323526
$ws.Range("E11").Value = 15.9457770229039
$ws.Range("G11").Value = 3.678395645801857
$ws.Range("I11").Value = 25.30869031997323
$ws.Range("J11").Value = 9.205134022248968
$ws.Range("K11").Value = 12.04288067508407
$ws.Range("M11").Value = 17.67267878938313
$ws.Range("N11").Value = 20.59689207500574
$ws.Range("O11").Value = 28.52746084600884
# Row 12
$ws.Range("B12").Value = 11.8415268972568
$ws.Range("C12").Value = 7.141392908160552
$ws.Range("D12").Value = 14.89938283321135
$ws.Range("E12").Value = 15.95210610032186
$ws.Range("G12").Value = 3.677905703932997
$ws.Range("I12").Value = 25.29310171082483
$ws.Range("J12").Value = 9.201562573158988
$ws.Range("K12").Value = 12.08564703776934
$ws.Range("M12").Value = 17.69301978640468
$ws.Range("N12").Value = 20.58107823612577
$ws.Range("O12").Value = 28.5172624313935
# Row 13
$ws.Range("B13").Value = 11.82765849200858
$ws.Range("C13").Value = 7.13297063797035
$ws.Range("D13").Value = 14.89633547446198
$ws.Range("E13").Value = 15.9507243702212
$ws.Range("G13").Value = 3.678010799904951
$ws.Range("I13").Value = 25.29643062425275
$ws.Range("J13").Value = 9.202325888216279
$ws.Range("K13").Value = 12.07643660369118
$ws.Range("M13").Value = 17.68862184563243
$ws.Range("N13").Value = 20.58447251996415
$ws.Range("O13").Value = 28.51942343775868
# Row 14
$ws.Range("B14").Value = 11.78236322361724
$ws.Range("C14").Value = 7.105463890222002
$ws.Range("D14").Value = 14.88646610353049
$ws.Range("E14").Value = 15.94628931073093
$ws.Range("G14").Value = 3.678355147860677
$ws.Range("I14").Value = 25.30739505177108
$ws.Range("J14").Value = 9.204837555725105
$ws.Range("K14").Value = 12.04639906626149
$ws.Range("M14").Value = 17.67434451292008
$ws.Range("N14").Value = 20.59558587297139
$ws.Range("O14").Value = 28.52660586355725
# Row 15
$ws.Range("B15").Value = 11.75456789933381
$ws.Range("C15").Value = 7.088585229846424
$ws.Range("D15").Value = 14.88047378876439
$ws.Range("E15").Value = 15.94362737788894
$ws.Range("G15").Value = 3.678567306735232
$ws.Range("I15").Value = 25.31419414355547
$ws.Range("J15").Value = 9.206393189750534
$ws.Range("K15").Value = 12.02800063752552
$ws.Range("M15").Value = 17.66564963618427
$ws.Range("N15").Value = 20.60242684277285
$ws.Range("O15").Value = 28.53110896047985
# Row 16
$ws.Range("B16").Value = 11.59455437846076
$ws.Range("C16").Value = 6.991429099987374
$ws.Range("D16").Value = 14.84695198625894
$ws.Range("E16").Value = 15.9292133906714
$ws.Range("G16").Value = 3.679802080044858
$ws.Range("I16").Value = 25.35442945925737
$ws.Range("J16").Value = 9.215571218921674
$ws.Range("K16").Value = 11.92260102391969
$ws.Range("M16").Value = 17.61660423186033
$ws.Range("N16").Value = 20.64214737005008
$ws.Range("O16").Value = 28.55850238108416
# Row 17
$ws.Range("B17").Value = 11.49579209332221
$ws.Range("C17").Value = 6.93147340323343
$ws.Range("D17").Value = 14.82712349118994
$ws.Range("E17").Value = 15.9211235116958
$ws.Range("G17").Value = 3.680576544757186
$ws.Range("I17").Value = 25.38024723710289
$ws.Range("J17").Value = 9.221436850729262
$ws.Range("K17").Value = 11.85800279946569
$ws.Range("M17").Value = 17.58722326153373
$ws.Range("N17").Value = 20.66697750957139
$ws.Range("O17").Value = 28.57672349392385
# Row 18
$ws.Range("B18").Value = 11.43878752233678
$ws.Range("C18").Value = 6.896871244390998
$ws.Range("D18").Value = 14.81599161868997
$ws.Range("E18").Value = 15.9167480324867
$ws.Range("G18").Value = 3.681028248428774
$ws.Range("I18").Value = 25.39551364958826
$ws.Range("J18").Value = 9.224897141491088
$ws.Range("K18").Value = 11.82088215844858
$ws.Range("M18").Value = 17.5705875601764
$ws.Range("N18").Value = 20.68142962831726
$ws.Range("O18").Value = 28.58772385764421
# Row 19
$ws.Range("B19").Value = 11.41945459722866
$ws.Range("C19").Value = 6.885136644351367
$ws.Range("D19").Value = 14.81226966495821
$ws.Range("E19").Value = 15.91531434158266
$ws.Range("G19").Value = 3.681182262768833
$ws.Range("I19").Value = 25.40075413658258
$ws.Range("J19").Value = 9.226083607982423
$ws.Range("K19").Value = 11.80832083747136
$ws.Range("M19").Value = 17.56500059199192
$ws.Range("N19").Value = 20.68635218274589
$ws.Range("O19").Value = 28.59153767622191
# Row 20
$ws.Range("B20").Value = 11.5063265976864
$ws.Range("C20").Value = 6.937868213037603
$ws.Range("D20").Value = 14.82920607305788
$ws.Range("E20").Value = 15.92195598420468
$ws.Range("G20").Value = 3.680493454945054
$ws.Range("I20").Value = 25.37745575599091
$ws.Range("J20").Value = 9.220803490381863
$ws.Range("K20").Value = 11.8648760895445
$ws.Range("M20").Value = 17.59032373188511
$ws.Range("N20").Value = 20.66431666325444
$ws.Range("O20").Value = 28.57472999233349
# Row 21
$ws.Range("B21").Value = 11.7956782108856
$ws.Range("C21").Value = 7.113549603900241
$ws.Range("D21").Value = 14.88935395477656
$ws.Range("E21").Value = 15.94758060843056
$ws.Range("G21").Value = 3.678253747060921
$ws.Range("I21").Value = 25.30415722171804
$ws.Range("J21").Value = 9.204096241429136
$ws.Range("K21").Value = 12.05522179559163
$ws.Range("M21").Value = 17.67852762942654
$ws.Range("N21").Value = 20.59231458658747
$ws.Range("O21").Value = 28.5244746079581
# Row 22
$ws.Range("B22").Value = 11.98257455888454
$ws.Range("C22").Value = 7.227059308316973
$ws.Range("D22").Value = 14.93104116232994
$ws.Range("E22").Value = 15.96677744325813
$ws.Range("G22").Value = 3.676845321076147
$ws.Range("I22").Value = 25.25996932435109
$ws.Range("J22").Value = 9.193945619309519
$ws.Range("K22").Value = 12.17967549328548
$ws.Range("M22").Value = 17.73843993749352
$ws.Range("N22").Value = 20.54676730041547
$ws.Range("O22").Value = 28.49626816340663
# Row 23
$ws.Range("B23").Value = 11.88304729205734
$ws.Range("C23").Value = 7.166609083962743
$ws.Range("D23").Value = 14.90857710604522
$ws.Range("E23").Value = 15.95630873844797
$ws.Range("G23").Value = 3.677591975237402
$ws.Range("I23").Value = 25.28321284393868
$ws.Range("J23").Value = 9.199292976195585
$ws.Range("K23").Value = 12.11325971587443
$ws.Range("M23").Value = 17.70626025283234
$ws.Range("N23").Value = 20.57093893630454
$ws.Range("O23").Value = 28.51089775049242
# Row 24
$ws.Range("B24").Value = 11.50156464352786
$ws.Range("C24").Value = 6.934977530250974
$ws.Range("D24").Value = 14.82826370307467
$ws.Range("E24").Value = 15.92157876499574
$ws.Range("G24").Value = 3.680530999743237
$ws.Range("I24").Value = 25.37871646569392
$ws.Range("J24").Value = 9.221089558268988
$ws.Range("K24").Value = 11.86176861695924
$ws.Range("M24").Value = 17.58892121153596
$ws.Range("N24").Value = 20.66551908060568
$ws.Range("O24").Value = 28.57562961973041
# Row 25
$ws.Range("B25").Value = 11.08275631221313
$ws.Range("C25").Value = 6.680821286693933
$ws.Range("D25").Value = 14.75209411726344
$ws.Range("E25").Value = 15.89473979440726
$ws.Range("G25").Value = 3.683940443767583
$ws.Range("I25").Value = 25.49766036363654
$ws.Range("J25").Value = 9.247910315052914
$ws.Range("K25").Value = 11.59197183895828
$ws.Range("M25").Value = 17.47246597872063
$ws.Range("N25").Value = 20.77406187925678
$ws.Range("O25").Value = 28.66533014962067
